$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (new Police Commissioner, volume/issue number, report week dates) ---
$ws.Range("M6").Value = "Jessica S. Tisch"
$ws.Range("A8").Value = "Volume 31   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# --- Cells whose style/type changes (use Copy to replicate exact style index) ---
$ws.Range("I14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("I14").Copy($ws.Range("G15"))
$ws.Range("G15").Value = 1
$ws.Range("K14").Copy($ws.Range("H15"))
$ws.Range("H15").Value = -100
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("I14").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 2
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("C14").Copy($ws.Range("F27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))

# --- Cells whose value changes but style/type remains the same ---
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 41.666666666666
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -83.333333333333
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 178
$ws.Range("J16").Value = 182
$ws.Range("K16").Value = -2.197802197802
$ws.Range("L16").Value = 8.536585365853
$ws.Range("M16").Value = -48.554913294797
$ws.Range("N16").Value = -85.092127303182
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 284
$ws.Range("J17").Value = 288
$ws.Range("K17").Value = -1.388888888888
$ws.Range("L17").Value = 10.077519379845
$ws.Range("M17").Value = 49.473684210526
$ws.Range("N17").Value = -48.269581056466
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 284
$ws.Range("J18").Value = 263
$ws.Range("K18").Value = 7.984790874524
$ws.Range("L18").Value = 0.353356890459
$ws.Range("M18").Value = -35.891647855530
$ws.Range("N18").Value = -77.207062600321
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 27
$ws.Range("E19").Value = -59.259259259259
$ws.Range("F19").Value = 43
$ws.Range("H19").Value = -37.681159420289
$ws.Range("I19").Value = 651
$ws.Range("J19").Value = 727
$ws.Range("K19").Value = -10.453920220082
$ws.Range("L19").Value = 4.830917874396
$ws.Range("M19").Value = 40.301724137931
$ws.Range("N19").Value = 27.397260273972
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -55.555555555555
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 114
$ws.Range("J20").Value = 178
$ws.Range("K20").Value = -35.955056179775
$ws.Range("L20").Value = -32.142857142857
$ws.Range("M20").Value = -22.972972972973
$ws.Range("N20").Value = -86.851211072664
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = -53.125
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 165
$ws.Range("H21").Value = -32.121212121212
$ws.Range("I21").Value = 1531
$ws.Range("J21").Value = 1655
$ws.Range("K21").Value = -7.492447129909
$ws.Range("L21").Value = 1.457919151756
$ws.Range("M21").Value = -4.372267332916
$ws.Range("N21").Value = -65.369825831259
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = -12.5
$ws.Range("L22").Value = -44.736842105263
$ws.Range("M22").Value = -48.780487804878
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -75
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = -80
$ws.Range("I23").Value = 150
$ws.Range("J23").Value = 190
$ws.Range("K23").Value = -21.052631578947
$ws.Range("L23").Value = 1.351351351351
$ws.Range("M23").Value = 19.047619047619
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -37.5
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 999
$ws.Range("J24").Value = 951
$ws.Range("K24").Value = 5.047318611987
$ws.Range("L24").Value = -8.933454876937
$ws.Range("M24").Value = -13.506493506493
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -80
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 296
$ws.Range("J25").Value = 146
$ws.Range("K25").Value = 102.739726027397
$ws.Range("L25").Value = 19.354838709677
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -37.5
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 477
$ws.Range("J26").Value = 452
$ws.Range("K26").Value = 5.530973451327
$ws.Range("L26").Value = 2.580645161290
$ws.Range("M26").Value = 3.246753246753
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 69
$ws.Range("K28").Value = 38
$ws.Range("L28").Value = 53.333333333333
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("L33").Value = 25
